# Purchase Requisition to Post Payment
# ------------------------------------------------------------------
# 1. POReceipt: update the QuantityBefore/QuantityAfter result values
#    written out by the latest automation run (row 2, columns O/P).
#    The NumberFormat "@" / ClearFormats dance forces Excel to store
#    the value as literal text ("1066.0") instead of silently
#    re-interpreting the numeric-looking string as a number, while
#    leaving the cell's style untouched (matching the source data,
#    which carries no explicit style on these cells).
# 2. POIssue: the user had clicked into D2 on this sheet at some
#    point during the session.
# 3. AllocatePOComponent: a new automation run id was written into
#    H2, and this sheet is where the user ended up / left the
#    selection when they saved.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- POReceipt: refresh the latest QuantityBefore/QuantityAfter run results ---
$wsReceipt = $wb.Worksheets.Item("POReceipt")

$wsReceipt.Range("O2").NumberFormat = "@"
$wsReceipt.Range("O2").Value = "1066.0"
$wsReceipt.Range("O2").ClearFormats()

$wsReceipt.Range("P2").NumberFormat = "@"
$wsReceipt.Range("P2").Value = "1076.0"
$wsReceipt.Range("P2").ClearFormats()

# --- POIssue: leave the cursor on D2 as the user left it ---
$wsIssue = $wb.Worksheets.Item("POIssue")
$wsIssue.Select()
$wsIssue.Range("D2").Select()

# --- AllocatePOComponent: new record id + final active selection ---
$wsAllocate = $wb.Worksheets.Item("AllocatePOComponent")
$wsAllocate.Select()
$wsAllocate.Range("H2").Value = "a2S1K000001xELw"
$wsAllocate.Range("H2").Select()
